$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 5 (pushes the existing row 5..111 down to
# 6..112, carrying their values/formatting with them and growing the
# sheet dimension to A1:R112 automatically).
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new price record.
$ws.Range("A5").Value = 5
$ws.Range("B5").Value = "Macroferia Regional de Talca"
$ws.Range("C5").Value = "Maule"
$ws.Range("D5").Value = 44882
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 100112026
$ws.Range("G5").Value = "Haba"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 9000
$ws.Range("L5").Value = 9000
$ws.Range("M5").Value = 9000
$ws.Range("N5").Value = "$/saco 25 kilos"
$ws.Range("O5").Value = "Región del Maule"
$ws.Range("P5").Value = 360
$ws.Range("Q5").Value = 25
$ws.Range("R5").Value = "Hortaliza"
